$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5700.2
$ws.Range("I40").Value = 6835.1
$ws.Range("J40").Value = 3430.4
$ws.Range("K40").Value = 6835.1
$ws.Range("L40").Value = 3430.4
$ws.Range("M40").Value = -6660.1
$ws.Range("N40").Value = -3780.4
$ws.Range("H51").Value = 10106038
$ws.Range("I51").Value = 27780228
$ws.Range("J51").Value = 6500.2856
$ws.Range("K51").Value = 27780228
$ws.Range("L51").Value = 6500.2856
$ws.Range("M51").Value = -27779744
$ws.Range("N51").Value = -7468.2856
$ws.Range("H74").Value = 4021.3572
$ws.Range("J74").Value = 4599.8
$ws.Range("L74").Value = 4599.8
$ws.Range("N74").Value = -6471.8
$ws.Range("H77").Value = 4021.3572
$ws.Range("J77").Value = 4599.8
$ws.Range("L77").Value = 22999
$ws.Range("N77").Value = -32359
$ws.Range("H109").Value = 38508
$ws.Range("J109").Value = 38508
$ws.Range("L109").Value = 38508
$ws.Range("N109").Value = -41282
$ws.Range("H112").Value = 2025.2632
$ws.Range("J112").Value = 2025.2632
$ws.Range("L112").Value = 6075.7896
$ws.Range("N112").Value = -8291.7896
$ws.Range("H128").Value = 50766.668
$ws.Range("J128").Value = 50766.668
$ws.Range("L128").Value = 50766.668
$ws.Range("N128").Value = -60726.668
$ws.Range("H129").Value = 1162.519
$ws.Range("J129").Value = 1118.6232
$ws.Range("L129").Value = 3355.8696
$ws.Range("N129").Value = -13355.8696
$ws.Range("H130").Value = 44995.555
$ws.Range("J130").Value = 44995.555
$ws.Range("L130").Value = 44995.555
$ws.Range("N130").Value = -55035.555
$ws.Range("H133").Value = 53795.625
$ws.Range("J133").Value = 53795.625
$ws.Range("L133").Value = 53795.625
$ws.Range("N133").Value = -63915.625
$ws.Range("H138").Value = 1813.7526
$ws.Range("I138").Value = 1462.5834
$ws.Range("J138").Value = 2157.7551
$ws.Range("K138").Value = 4387.7502
$ws.Range("L138").Value = 6473.265299999999
$ws.Range("M138").Value = 752.2497999999996
$ws.Range("N138").Value = -16753.2653

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3428.9768
$ws.Range("I2").Value = 3790.658
$ws.Range("J2").Value = 680.2
$ws.Range("K2").Value = 3790.658
$ws.Range("L2").Value = 680.2
$ws.Range("M2").Value = -3677.658
$ws.Range("N2").Value = -906.2
$ws.Range("H80").Value = 39997.332
$ws.Range("J80").Value = 49946
$ws.Range("L80").Value = 49946
$ws.Range("N80").Value = -51942
$ws.Range("H83").Value = 39997.332
$ws.Range("J83").Value = 49946
$ws.Range("L83").Value = 149838
$ws.Range("N83").Value = -159822
$ws.Range("H113").Value = 39078
$ws.Range("J113").Value = 39078
$ws.Range("L113").Value = 39078
$ws.Range("N113").Value = -47756
$ws.Range("H116").Value = 3428.9768
$ws.Range("I116").Value = 3790.658
$ws.Range("J116").Value = 680.2
$ws.Range("K116").Value = 3790.658
$ws.Range("L116").Value = 680.2
$ws.Range("M116").Value = -1496.658
$ws.Range("N116").Value = -5268.2
$ws.Range("H130").Value = 37975.8
$ws.Range("J130").Value = 37975.8
$ws.Range("L130").Value = 37975.8
$ws.Range("N130").Value = -48015.8
$ws.Range("H131").Value = 50563.668
$ws.Range("J131").Value = 50563.668
$ws.Range("L131").Value = 50563.668
$ws.Range("N131").Value = -60643.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3428.9768
$ws.Range("I3").Value = 3790.658
$ws.Range("J3").Value = 680.2
$ws.Range("K3").Value = 3790.658
$ws.Range("L3").Value = 680.2
$ws.Range("M3").Value = -3676.658
$ws.Range("N3").Value = -908.2
$ws.Range("H86").Value = 2083.8462
$ws.Range("I86").Value = 2221.111
$ws.Range("J86").Value = 1775
$ws.Range("K86").Value = 2221.111
$ws.Range("L86").Value = 1775
$ws.Range("M86").Value = -1098.111
$ws.Range("N86").Value = -4021
$ws.Range("H89").Value = 2083.8462
$ws.Range("I89").Value = 2221.111
$ws.Range("J89").Value = 1775
$ws.Range("K89").Value = 11105.555
$ws.Range("L89").Value = 8875
$ws.Range("M89").Value = -5489.555
$ws.Range("N89").Value = -20107
$ws.Range("H126").Value = 50772
$ws.Range("J126").Value = 50772
$ws.Range("L126").Value = 50772
$ws.Range("N126").Value = -60652
$ws.Range("H130").Value = 48815
$ws.Range("J130").Value = 48815
$ws.Range("L130").Value = 48815
$ws.Range("N130").Value = -58855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 44887.75
$ws.Range("J20").Value = 44887.75
$ws.Range("L20").Value = 44887.75
$ws.Range("N20").Value = -45359.75
$ws.Range("H30").Value = 44887.75
$ws.Range("J30").Value = 44887.75
$ws.Range("L30").Value = 44887.75
$ws.Range("N30").Value = -45069.75
$ws.Range("H100").Value = 46996
$ws.Range("J100").Value = 46996
$ws.Range("L100").Value = 46996
$ws.Range("N100").Value = -49160
$ws.Range("H116").Value = 47822.332
$ws.Range("J116").Value = 47822.332
$ws.Range("L116").Value = 47822.332
$ws.Range("N116").Value = -57000.332
$ws.Range("H128").Value = 44887.75
$ws.Range("J128").Value = 44887.75
$ws.Range("L128").Value = 44887.75
$ws.Range("N128").Value = -54847.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4416.5
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 5199.8
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 15599.4
$ws.Range("M70").Value = -1185
$ws.Range("N70").Value = -16229.4
$ws.Range("H73").Value = 4416.5
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 5199.8
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 15599.4
$ws.Range("M73").Value = -408
$ws.Range("N73").Value = -17783.4
$ws.Range("H92").Value = 1256.375
$ws.Range("I92").Value = 1099.6364
$ws.Range("J92").Value = 1601.2
$ws.Range("K92").Value = 3298.9092
$ws.Range("L92").Value = 4803.6
$ws.Range("M92").Value = -2050.9092
$ws.Range("N92").Value = -7299.6
$ws.Range("H133").Value = 9609.166999999999
$ws.Range("J133").Value = 6485
$ws.Range("L133").Value = 19455
$ws.Range("N133").Value = -29575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 20000
$ws.Range("J40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("N40").Value = -20302
$ws.Range("H113").Value = 8363.9375
$ws.Range("I113").Value = 11055.454
$ws.Range("K113").Value = 11055.454
$ws.Range("M113").Value = -8885.454
$ws.Range("H118").Value = 37204
$ws.Range("J118").Value = 37204
$ws.Range("L118").Value = 37204
$ws.Range("N118").Value = -40518
$ws.Range("H130").Value = 52860.8
$ws.Range("J130").Value = 52860.8
$ws.Range("L130").Value = 52860.8
$ws.Range("N130").Value = -62900.8
$ws.Range("H136").Value = 14368.258
$ws.Range("J136").Value = 14368.258
$ws.Range("L136").Value = 43104.774
$ws.Range("N136").Value = -48204.774

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3164.625
$ws.Range("I7").Value = 2522.4
$ws.Range("K7").Value = 2522.4
$ws.Range("M7").Value = -2410.4
$ws.Range("H46").Value = 2790.1667
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -812
$ws.Range("H126").Value = 3164.625
$ws.Range("I126").Value = 2522.4
$ws.Range("K126").Value = 7567.200000000001
$ws.Range("M126").Value = -5097.200000000001
$ws.Range("H130").Value = 48318
$ws.Range("J130").Value = 48318
$ws.Range("L130").Value = 48318
$ws.Range("N130").Value = -58358
$ws.Range("H134").Value = 44398
$ws.Range("J134").Value = 44398
$ws.Range("L134").Value = 44398
$ws.Range("N134").Value = -54538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H46").Value = 74800.17999999999
$ws.Range("J46").Value = 74800.17999999999
$ws.Range("L46").Value = 74800.17999999999
$ws.Range("N46").Value = -75262.17999999999
$ws.Range("H117").Value = 46052.25
$ws.Range("J117").Value = 46052.25
$ws.Range("L117").Value = 46052.25
$ws.Range("N117").Value = -55230.25
$ws.Range("H134").Value = 74800.17999999999
$ws.Range("J134").Value = 74800.17999999999
$ws.Range("L134").Value = 224400.54
$ws.Range("N134").Value = -229470.54
